# Generate Report for Handoff
# Adds a new row (row 3) to each of the three tables (Overview, zh-cn, de-de)
# describing the newly-handed-off file
# "3b4ac6e7-2f6f-4084-a987-7f5131021839.md".

$wb = $excel.ActiveWorkbook

$guidFile       = "3b4ac6e7-2f6f-4084-a987-7f5131021839.md"
$guidPath       = "e2e\3b4ac6e7-2f6f-4084-a987-7f5131021839.md"
$baseUrl        = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dda0888eddfd8ef070686467865b3726c451d23a/e2e/"
$dateTimeFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item("Overview")
$lo.ListRows.Add() | Out-Null

$ws.Hyperlinks.Add($ws.Range("B3"), ($baseUrl + $guidFile), $null, $null, $guidPath) | Out-Null

$ws.Range("A3").Value = "'" + $guidFile
$ws.Range("C3").Value = "'.md"
$ws.Range("D3").Value = "'"
$ws.Range("E3").Value = "'Ready for handoff"
$ws.Range("F3").Value = "'Ready for handoff"
$ws.Range("G3").Value = "2016-08-22 20:39:30"
$ws.Range("G3").NumberFormat = $dateTimeFormat

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item("zh-cn")
$lo.ListRows.Add() | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), ($baseUrl + $guidFile), $null, $null, $guidFile) | Out-Null

$ws.Range("B3").Value = "'.md"
$ws.Range("C3").Value = "'Ready for handoff"
$ws.Range("D3").Value = "'e2e"
$ws.Range("E3").Value = "'ht"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "'3b4ac6e7-2f6f-4084-a987-7f5131021839.ecfbbcb7048f8754fcfbc3ab24c1bcf4c23ce337.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-22 20:39:26"
$ws.Range("H3").NumberFormat = $dateTimeFormat
$ws.Range("I3").Value = "'"
$ws.Range("J3").Value = "'"
$ws.Range("K3").Value = "0001-01-01 00:00:00"
$ws.Range("K3").NumberFormat = $dateTimeFormat
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = "'"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item("de-de")
$lo.ListRows.Add() | Out-Null

$ws.Hyperlinks.Add($ws.Range("A3"), ($baseUrl + $guidFile), $null, $null, $guidFile) | Out-Null

$ws.Range("B3").Value = "'.md"
$ws.Range("C3").Value = "'Ready for handoff"
$ws.Range("D3").Value = "'e2e"
$ws.Range("E3").Value = "'ht"
$ws.Range("F3").Value = "'False"
$ws.Range("G3").Value = "'3b4ac6e7-2f6f-4084-a987-7f5131021839.ecfbbcb7048f8754fcfbc3ab24c1bcf4c23ce337.de-de.xlf"
$ws.Range("H3").Value = "2016-08-22 20:39:30"
$ws.Range("H3").NumberFormat = $dateTimeFormat
$ws.Range("I3").Value = "'"
$ws.Range("J3").Value = "'"
$ws.Range("K3").Value = "0001-01-01 00:00:00"
$ws.Range("K3").NumberFormat = $dateTimeFormat
$ws.Range("L3").Value = "'"
$ws.Range("M3").Value = "'True"
$ws.Range("N3").Value = "'"
$ws.Range("O3").Value = "'False"
$ws.Range("P3").Value = "'"
